# Apply the "Updated cryptos list" data refresh to Sheet1.
# Only cells whose values actually changed are touched; row/column
# layout, headers and styles are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '37.384.66'
$ws.Range("E2").Value = '  +1.87%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '2.078.59'
$ws.Range("E3").Value = '  -1.73%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.04%  '
# Row 5: BNB
$ws.Range("D5").Value = '''251.05'
$ws.Range("E5").Value = '  +0.39%  '
# Row 6: XRP
$ws.Range("E6").Value = '  -1.21%  '
# Row 7: USDC
$ws.Range("E7").Value = '  -0.02%  '
# Row 8: Solana
$ws.Range("D8").Value = '''57.05'
$ws.Range("E8").Value = '  +25.84%  '
# Row 9: OKB
$ws.Range("D9").Value = '''62.22'
$ws.Range("E9").Value = '  +1.68%  '
# Row 10: Cardano
$ws.Range("D10").Value = '''0.387'
$ws.Range("E10").Value = '  +5.20%  '
# Row 11: Dogecoin
$ws.Range("D11").Value = '''0.0753'
$ws.Range("E11").Value = '  +2.88%  '
# Row 13: Chainlink
$ws.Range("D13").Value = '''15.64'
$ws.Range("E13").Value = '  +6.94%  '
# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '2.384.66'
$ws.Range("E14").Value = '  -1.34%  '
# Row 15: Polygon
$ws.Range("D15").Value = '''0.844'
$ws.Range("E15").Value = '  -0.36%  '
# Row 16: Polkadot
$ws.Range("D16").Value = '''5.31'
$ws.Range("E16").Value = '  +4.83%  '
# Row 17: WrappedEther
$ws.Range("D17").Value = '2.081.20'
$ws.Range("E17").Value = '  -1.52%  '
# Row 18: WrappedBTC
$ws.Range("D18").Value = '37.278.76'
$ws.Range("E18").Value = '  +1.59%  '
# Row 19: Litecoin
$ws.Range("D19").Value = '''73.25'
$ws.Range("E19").Value = '  +0.11%  '
# Row 20: Avalanche
$ws.Range("D20").Value = '''14.83'
$ws.Range("E20").Value = '  +14.52%  '
# Row 21: ShibaInu
$ws.Range("D21").Value = '0.0₃0850'
$ws.Range("E21").Value = '  +3.43%  '
# Row 22: BitcoinCash
$ws.Range("D22").Value = '''240.88'
$ws.Range("E22").Value = '  -0.11%  '
# Row 23: Uniswap
$ws.Range("D23").Value = '''5.30'
$ws.Range("E23").Value = '  +4.14%  '
# Row 24: Dai
$ws.Range("E24").Value = '  -0.13%  '
# Row 25: Toncoin
$ws.Range("E25").Value = '  +0.25%  '
# Row 26: Monero
$ws.Range("D26").Value = '''171.82'
$ws.Range("E26").Value = '  +1.00%  '
# Row 27: Cosmos
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '''9.25'
$ws.Range("E27").Value = '  +2.38%  '
# Row 28: EthereumClassic
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''21.09'
$ws.Range("E28").Value = '  +2.31%  '
# Row 29: PancakeSwap
$ws.Range("E29").Value = '  +0.48%  '
# Row 30: Stellar
$ws.Range("E30").Value = '  +0.62%  '
# Row 31: Gas
$ws.Range("D31").Value = '''23.50'
$ws.Range("E31").Value = '  +5.43%  '
# Row 32: ImmutableX
$ws.Range("D32").Value = '''1.11'
$ws.Range("E32").Value = '  +22.21%  '
# Row 33: Filecoin
$ws.Range("D33").Value = '''4.58'
$ws.Range("E33").Value = '  +3.11%  '
# Row 34: Hedera
$ws.Range("D34").Value = '''0.0630'
$ws.Range("E34").Value = '  +5.46%  '
# Row 35: InternetComputer(DFINITY)
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = '''4.34'
$ws.Range("E35").Value = '  +6.28%  '
# Row 36: Kaspa
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.0906'
$ws.Range("E36").Value = '  -0.14%  '
# Row 37: BinanceUSD
$ws.Range("E37").Value = '  -0.02%  '
# Row 38: WEMIXToken
$ws.Range("E38").Value = '  -1.51%  '
# Row 39: LidoDAOToken
$ws.Range("D39").Value = '''2.29'
$ws.Range("E39").Value = '  -1.43%  '
# Row 40: TrustWalletToken
$ws.Range("E40").Value = '  -0.61%  '
# Row 41: VeChain
$ws.Range("E41").Value = '  +4.41%  '
# Row 42: Cronos
$ws.Range("D42").Value = '''0.100'
$ws.Range("E42").Value = '  +19.86%  '
# Row 43: InjectiveProtocol
$ws.Range("D43").Value = '''17.78'
$ws.Range("E43").Value = '  +8.94%  '
# Row 44: ARBITRUM
$ws.Range("E44").Value = '  -2.14%  '
# Row 45: Aave
$ws.Range("D45").Value = '''100.14'
$ws.Range("E45").Value = '  +0.14%  '
# Row 46: HuobiToken
$ws.Range("B46").Value = 'HuobiToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D46").Value = '''2.80'
$ws.Range("E46").Value = '  +0.06%  '
# Row 47: FTXToken
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '''4.20'
$ws.Range("E47").Value = '  +94.22%  '
# Row 48: Maker
$ws.Range("D48").Value = '1.326.92'
$ws.Range("E48").Value = '  -2.52%  '
# Row 49: RenderToken
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''2.38'
$ws.Range("E49").Value = '  +4.15%  '
# Row 50: MXToken
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '''2.94'
$ws.Range("E50").Value = '  +2.82%  '
# Row 51: FraxShare
$ws.Range("E51").Value = '  +8.61%  '
